$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.535955190658569
$ws.Range("B1").Value = 2.746290683746338
$ws.Range("C1").Value = 2.011159420013428
$ws.Range("D1").Value = 1.855700850486755
$ws.Range("E1").Value = 1.715347766876221
